$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "48.060.42"
Set-TextValue "E2" "  -0.21%  "
Set-TextValue "D3" "2.492.01"
Set-TextValue "E3" "  -1.39%  "
Set-TextValue "E4" "  -0.12%  "
Set-TextValue "D5" "317.64"
Set-TextValue "E5" "  -1.96%  "
Set-TextValue "D6" "105.53"
Set-TextValue "E6" "  -3.17%  "
Set-TextValue "E7" "  -1.88%  "
Set-TextValue "E8" "  -0.03%  "
Set-TextValue "D9" "0.537"
Set-TextValue "E9" "  -3.58%  "
Set-TextValue "D10" "38.90"
Set-TextValue "E10" "  -5.02%  "
Set-TextValue "D11" "20.21"
Set-TextValue "E11" "  -1.37%  "
Set-TextValue "E12" "  -3.02%  "
Set-TextValue "E13" "  +0.20%  "
Set-TextValue "E14" "  -2.66%  "
Set-TextValue "D15" "2.885.48"
Set-TextValue "E15" "  -1.41%  "
Set-TextValue "D16" "2.491.62"
Set-TextValue "E16" "  -1.45%  "
Set-TextValue "D17" "0.828"
Set-TextValue "E17" "  -3.64%  "
Set-TextValue "D18" "47.982.00"
Set-TextValue "E18" "  -0.08%  "
Set-TextValue "D19" "2.99"
Set-TextValue "E19" "  +10.73%  "
Set-TextValue "D20" "12.77"
Set-TextValue "E20" "  -3.91%  "
Set-TextValue "E21" "  -1.00%  "
Set-TextValue "E22" "  -2.43%  "
Set-TextValue "D23" "71.02"
Set-TextValue "E23" "  -1.68%  "
Set-TextValue "D24" "268.24"
Set-TextValue "E24" "  -0.54%  "
Set-TextValue "D25" "2.50"
Set-TextValue "E25" "  -3.13%  "
Set-TextValue "E26" "  +0.18%  "
Set-TextValue "D27" "25.71"
Set-TextValue "E27" "  -1.81%  "
Set-TextValue "D28" "2.28"
Set-TextValue "E28" "  +3.11%  "
Set-TextValue "D29" "9.73"
Set-TextValue "E29" "  -4.17%  "
Set-TextValue "D30" "0.139"
Set-TextValue "E30" "  -2.17%  "
Set-TextValue "D31" "34.51"
Set-TextValue "E31" "  -2.96%  "
Set-TextValue "D32" "49.23"
Set-TextValue "E32" "  -0.73%  "
Set-TextValue "E33" "  -0.12%  "
Set-TextValue "D34" "19.06"
Set-TextValue "D35" "5.27"
Set-TextValue "E35" "  -2.59%  "
Set-TextValue "E36" "  -2.89%  "
Set-TextValue "D37" "1.94"
Set-TextValue "E37" "  -2.40%  "
Set-TextValue "D38" "4.58"
Set-TextValue "E38" "  -3.52%  "
Set-TextValue "E39" "  -4.47%  "
Set-TextValue "D40" "123.35"
Set-TextValue "E40" "  +3.38%  "
Set-TextValue "B41" "Stellar"
Set-TextValue "C41" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D41" "0.110"
Set-TextValue "E41" "  -1.77%  "
Set-TextValue "B42" "EnergySwap"
Set-TextValue "C42" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D42" "22.23"
Set-TextValue "E42" "  +0.76%  "
Set-TextValue "E43" "  +1.24%  "
Set-TextValue "E44" "  +0.31%  "
Set-TextValue "D45" "1.999.35"
Set-TextValue "E45" "  -0.77%  "
Set-TextValue "E46" "  -0.13%  "
Set-TextValue "E47" "  +0.48%  "
Set-TextValue "D49" "8.92"
Set-TextValue "E49" "  -2.65%  "
Set-TextValue "D50" "5.19"
Set-TextValue "E50" "  -1.46%  "
Set-TextValue "D51" "78.53"
Set-TextValue "E51" "  -1.32%  "
